$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 188, shifting all rows below it up by one.
$ws.Rows.Item(188).Delete()
